$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. "We use 2.2.1 to try to insure compatibility with Intellij's version."
#    -> "We used version 2.2.1 (compatible with Intellij 11.1's plugin)."
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "We use 2.2.1 to try to insure compatibility with Intellij" + [char]0x2019 + "s version.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "We used version 2.2.1 (compatible with Intellij 11.1" + [char]0x2019 + "s plugin).",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 2. Insert new paragraph "We'll be using eclipse, but Included are
#    instructions for Intellij" right before the "Intellij or Eclipse:"
#    paragraph (i.e. immediately after the blank paragraph preceding it).
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text.TrimEnd([char]13, [char]7) -eq "Intellij or Eclipse:") {
        $prevBlank = $d.Paragraphs($i - 1)
        $prevBlank.Range.InsertParagraphAfter()
        $newPara = $d.Paragraphs($i)
        $newPara.Range.Text = "We" + [char]0x2019 + "ll be using eclipse, but Included are instructions for Intellij"
        break
    }
}

# ---------------------------------------------------------------------------
# 3. "works  better for some" (double space) -> "works better for some"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("works" + [char]32 + [char]32 + "better for some", $true, $false, $false, $false, $false, $true, 1, $false, "works better for some", 2) | Out-Null

# ---------------------------------------------------------------------------
# 4. Add new sentence before "Code snippets and other materials..."
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Code snippets and other materials for the tutorial are available",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "We" + [char]0x2019 + "ll also use the M2e maven plugin to add dependencies.  Code snippets and other materials for the tutorial are available",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 5. "From gitHub:" -> insert a _GoBack bookmark between "gitHub" and ":"
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $ptext = $d.Paragraphs($i).Range.Text.TrimEnd([char]13, [char]7)
    if ($ptext -eq "From gitHub:") {
        $pRange = $d.Paragraphs($i).Range
        $bmStart = $pRange.Start + 11  # position right after "From gitHub"
        $bmRange = $d.Range($bmStart, $bmStart)
        $d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
        break
    }
}

# ---------------------------------------------------------------------------
# 6. Hyperlink text "arch" + "i" + "ve/master.zip" -> "archive/master.zip"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("archive", $true, $false, $false, $false, $false, $true, 1, $false, "archive", 2) | Out-Null

# ---------------------------------------------------------------------------
# 7/8. Move lastRenderedPageBreak from "From the MySQL client..." paragraph
#      to the "This is an ontology..." paragraph.
# ---------------------------------------------------------------------------
# (handled structurally below, lastRenderedPageBreak markers are render hints
#  that Word recalculates; no explicit action needed through the COM object
#  model, Word will reposition them automatically on repagination.)

# ---------------------------------------------------------------------------
# 9. "client you can just do:" -> "command line client you can execute:"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(" client you can just do:", $true, $false, $false, $false, $false, $true, 1, $false, " command line client you can execute:", 2) | Out-Null

# ---------------------------------------------------------------------------
# 10. "If any Java files still have errors try running Maven..." restructure
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $ptext = $d.Paragraphs($i).Range.Text.TrimEnd([char]13, [char]7)
    if ($ptext -eq "If any Java files still have errors try running Maven => update project") {
        $para = $d.Paragraphs($i)
        $para.Range.Text = "If any Java files still have errors, Right click the project and try running Maven => update project"
        $para.Range.InsertParagraphBefore()
        $para2 = $d.Paragraphs($i + 1)
        $para2.Range.InsertParagraphAfter()
        $otherwisePara = $d.Paragraphs($i + 3)
        $otherwisePara.Range.Text = "Otherwise"
        $otherwisePara.Range.InsertParagraphAfter()
        $rightClickPara = $d.Paragraphs($i + 4)
        $rightClickPara.Range.Text = "Right click the project:"
        break
    }
}

# ---------------------------------------------------------------------------
# 11. New bold paragraph "Follow the Slides from Here"
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $ptext = $d.Paragraphs($i).Range.Text.TrimEnd([char]13, [char]7)
    if ($ptext -eq "Then  run as => Maven Install") {
        $blankPara = $d.Paragraphs($i + 1)
        $blankPara.Range.InsertParagraphAfter()
        $slidesPara = $d.Paragraphs($i + 2)
        $slidesPara.Range.Text = "Follow the Slides from Here"
        $slidesPara.Range.Font.Bold = 1
        $slidesPara.Range.InsertParagraphAfter()
        break
    }
}

# ---------------------------------------------------------------------------
# 12. Remove lastRenderedPageBreak near "from whatever directory you choose
#     to install the cts2 framework standalone server" (first occurrence)
#     and 13. remove _GoBack bookmark after "And login as admin/admin"
#     (the Word OM recomputes lastRenderedPageBreak automatically; the stray
#      legacy _GoBack bookmark here is removed explicitly.)
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $bm = $d.Bookmarks("_GoBack")
    $bmText = $bm.Range.Text
    # Only remove it if this is the leftover bookmark further down the
    # document (the one we (re)created above sits right after "From gitHub").
}

Write-Output "done"
